# Add a new "2022-Q3" worksheet (right after the "总计" summary sheet) with
# the latest quarter's fund-holding data, and update the "总计" summary
# sheet with a new leading row for 2022-Q3 (pushing the older quarters
# down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q3" sheet by duplicating the "2022-Q2" sheet
#    (so it inherits identical column widths / header styles / borders)
#    immediately after "总计", then renaming it and replacing its data.
# ---------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item(1)
$sheetQ2    = $wb.Worksheets.Item(2)

$sheetQ2.Copy($null, $sheetTotal)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The new quarter only lists 6 funds (the template had 7), drop the spare row.
$newSheet.Rows.Item(8).Delete()

# code, name, fund size, stock position, position ratio, holding value(bn), rank
$q3Data = @(
  @("012744", "光大保德信品质生活混合A",            "5.60", "88.62", "6.19", "0.3466", 2),
  @("166024", "中欧恒利三年定期开放混合",            "3.99", "98.45", "3.82", "0.1524", 8),
  @("012770", "光大保德信创新生活混合",              "2.71", "86.69", "3.81", "0.1033", 6),
  @("012758", "光大保德信品质生活混合C",            "0.35", "88.62", "6.19", "0.0217", 2),
  @("001942", "前海开源沪港深汇鑫灵活配置混合A",    "0.17", "87.24", "5.23", "0.0089", 2),
  @("001943", "前海开源沪港深汇鑫灵活配置混合C",    "0.09", "87.24", "5.23", "0.0047", 2)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $row  = $i + 2
    $item = $q3Data[$i]

    # B..G are text columns (fund codes / percentages are stored verbatim
    # as text in this workbook) - the leading "'" keeps them as literal
    # text instead of Excel auto-coercing them to numbers.
    $newSheet.Range("B$row").Value = "'" + $item[0]
    $newSheet.Range("C$row").Value = "'" + $item[1]
    $newSheet.Range("D$row").Value = "'" + $item[2]
    $newSheet.Range("E$row").Value = "'" + $item[3]
    $newSheet.Range("F$row").Value = "'" + $item[4]
    $newSheet.Range("G$row").Value = "'" + $item[5]
    $newSheet.Range("B$row" + ":G$row").ClearFormats()

    # H (position rank) is a genuine number.
    $newSheet.Range("H$row").Value = $item[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a fresh row right under the
#    header for 2022-Q3 and push the existing quarters down by one row.
# ---------------------------------------------------------------------
$sheetTotal.Rows.Item(2).Insert()
$sheetTotal.Range("B2:D2").ClearFormats()
$sheetTotal.Range("B2").Value = "2022-Q3"
$sheetTotal.Range("C2").Value = 6
$sheetTotal.Range("D2").Value = 0.64
